$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Reln-Lrp8 NATMI LR-pair table with the new TPM-based values.
# Rows 2-13 (ECs/FAPs/MuSCs as sending cluster) get updated metrics, and four
# new rows (14-17) are added for Resolving-Mac as a sending cluster, completing
# the full 4x4 matrix of sending/target cluster combinations.

# Row 2: ECs -> ECs
$abcd = New-Object 'object[,]' 1,4
$abcd[0,0] = "ECs"
$abcd[0,1] = "Reln"
$abcd[0,2] = "Lrp8"
$abcd[0,3] = "ECs"
$ws.Range("A2:D2").Value = $abcd

$data = New-Object 'object[,]' 1,16
$data[0,0] = 3
$data[0,1] = 1
$data[0,2] = 0.06455366666666666
$data[0,3] = 0.193661
$data[0,4] = 0.01357839286814829
$data[0,5] = 0.01357839286814829
$data[0,6] = 3
$data[0,7] = 1
$data[0,8] = 0.4216986666666666
$data[0,9] = 1.265096
$data[0,10] = 0.2516921781598699
$data[0,11] = 0.2516921781598699
$data[0,12] = 0.02722219516177777
$data[0,13] = 0.2449997564559999
$data[0,14] = 0.003417575276894686
$data[0,15] = 0.003417575276894686
$ws.Range("E2:T2").Value = $data

# Row 3: ECs -> FAPs
$abcd = New-Object 'object[,]' 1,4
$abcd[0,0] = "ECs"
$abcd[0,1] = "Reln"
$abcd[0,2] = "Lrp8"
$abcd[0,3] = "FAPs"
$ws.Range("A3:D3").Value = $abcd

$data = New-Object 'object[,]' 1,16
$data[0,0] = 3
$data[0,1] = 1
$data[0,2] = 0.06455366666666666
$data[0,3] = 0.193661
$data[0,4] = 0.01357839286814829
$data[0,5] = 0.01357839286814829
$data[0,6] = 3
$data[0,7] = 1
$data[0,8] = 0.437958
$data[0,9] = 1.313874
$data[0,10] = 0.2613966125002536
$data[0,11] = 0.2613966125002536
$data[0,12] = 0.028271794746
$data[0,13] = 0.254446152714
$data[0,14] = 0.003549345898931566
$data[0,15] = 0.003549345898931566
$ws.Range("E3:T3").Value = $data

# Row 4: ECs -> MuSCs
$abcd = New-Object 'object[,]' 1,4
$abcd[0,0] = "ECs"
$abcd[0,1] = "Reln"
$abcd[0,2] = "Lrp8"
$abcd[0,3] = "MuSCs"
$ws.Range("A4:D4").Value = $abcd

$data = New-Object 'object[,]' 1,16
$data[0,0] = 3
$data[0,1] = 1
$data[0,2] = 0.06455366666666666
$data[0,3] = 0.193661
$data[0,4] = 0.01357839286814829
$data[0,5] = 0.01357839286814829
$data[0,6] = 3
$data[0,7] = 1
$data[0,8] = 0.4328273333333333
$data[0,9] = 1.298482
$data[0,10] = 0.2583343579312433
$data[0,11] = 0.2583343579312433
$data[0,12] = 0.02794059140022222
$data[0,13] = 0.251465322602
$data[0,14] = 0.003507765403331261
$data[0,15] = 0.003507765403331262
$ws.Range("E4:T4").Value = $data

# Row 5: ECs -> Resolving-Mac
$abcd = New-Object 'object[,]' 1,4
$abcd[0,0] = "ECs"
$abcd[0,1] = "Reln"
$abcd[0,2] = "Lrp8"
$abcd[0,3] = "Resolving-Mac"
$ws.Range("A5:D5").Value = $abcd

$data = New-Object 'object[,]' 1,16
$data[0,0] = 3
$data[0,1] = 1
$data[0,2] = 0.06455366666666666
$data[0,3] = 0.193661
$data[0,4] = 0.01357839286814829
$data[0,5] = 0.01357839286814829
$data[0,6] = 3
$data[0,7] = 1
$data[0,8] = 0.38297
$data[0,9] = 1.14891
$data[0,10] = 0.2285768514086331
$data[0,11] = 0.2285768514086331
$data[0,12] = 0.02472211772333333
$data[0,13] = 0.22249905951
$data[0,14] = 0.003103706288990774
$data[0,15] = 0.003103706288990775
$ws.Range("E5:T5").Value = $data

# Row 6: FAPs -> ECs
$abcd = New-Object 'object[,]' 1,4
$abcd[0,0] = "FAPs"
$abcd[0,1] = "Reln"
$abcd[0,2] = "Lrp8"
$abcd[0,3] = "ECs"
$ws.Range("A6:D6").Value = $abcd

$data = New-Object 'object[,]' 1,16
$data[0,0] = 3
$data[0,1] = 1
$data[0,2] = 0.06084700000000001
$data[0,3] = 0.182541
$data[0,4] = 0.01279872257472933
$data[0,5] = 0.01279872257472933
$data[0,6] = 3
$data[0,7] = 1
$data[0,8] = 0.4216986666666666
$data[0,9] = 1.265096
$data[0,10] = 0.2516921781598699
$data[0,11] = 0.2516921781598699
$data[0,12] = 0.02565909877066666
$data[0,13] = 0.230931888936
$data[0,14] = 0.003221338362497524
$data[0,15] = 0.003221338362497523
$ws.Range("E6:T6").Value = $data

# Row 7: FAPs -> FAPs
$abcd = New-Object 'object[,]' 1,4
$abcd[0,0] = "FAPs"
$abcd[0,1] = "Reln"
$abcd[0,2] = "Lrp8"
$abcd[0,3] = "FAPs"
$ws.Range("A7:D7").Value = $abcd

$data = New-Object 'object[,]' 1,16
$data[0,0] = 3
$data[0,1] = 1
$data[0,2] = 0.06084700000000001
$data[0,3] = 0.182541
$data[0,4] = 0.01279872257472933
$data[0,5] = 0.01279872257472933
$data[0,6] = 3
$data[0,7] = 1
$data[0,8] = 0.437958
$data[0,9] = 1.313874
$data[0,10] = 0.2613966125002536
$data[0,11] = 0.2613966125002536
$data[0,12] = 0.026648430426
$data[0,13] = 0.239835873834
$data[0,14] = 0.003345542725364772
$data[0,15] = 0.003345542725364771
$ws.Range("E7:T7").Value = $data

# Row 8: FAPs -> MuSCs
$abcd = New-Object 'object[,]' 1,4
$abcd[0,0] = "FAPs"
$abcd[0,1] = "Reln"
$abcd[0,2] = "Lrp8"
$abcd[0,3] = "MuSCs"
$ws.Range("A8:D8").Value = $abcd

$data = New-Object 'object[,]' 1,16
$data[0,0] = 3
$data[0,1] = 1
$data[0,2] = 0.06084700000000001
$data[0,3] = 0.182541
$data[0,4] = 0.01279872257472933
$data[0,5] = 0.01279872257472933
$data[0,6] = 3
$data[0,7] = 1
$data[0,8] = 0.4328273333333333
$data[0,9] = 1.298482
$data[0,10] = 0.2583343579312433
$data[0,11] = 0.2583343579312433
$data[0,12] = 0.02633624475133333
$data[0,13] = 0.237026202762
$data[0,14] = 0.003306349778682811
$data[0,15] = 0.003306349778682811
$ws.Range("E8:T8").Value = $data

# Row 9: FAPs -> Resolving-Mac
$abcd = New-Object 'object[,]' 1,4
$abcd[0,0] = "FAPs"
$abcd[0,1] = "Reln"
$abcd[0,2] = "Lrp8"
$abcd[0,3] = "Resolving-Mac"
$ws.Range("A9:D9").Value = $abcd

$data = New-Object 'object[,]' 1,16
$data[0,0] = 3
$data[0,1] = 1
$data[0,2] = 0.06084700000000001
$data[0,3] = 0.182541
$data[0,4] = 0.01279872257472933
$data[0,5] = 0.01279872257472933
$data[0,6] = 3
$data[0,7] = 1
$data[0,8] = 0.38297
$data[0,9] = 1.14891
$data[0,10] = 0.2285768514086331
$data[0,11] = 0.2285768514086331
$data[0,12] = 0.02330257559
$data[0,13] = 0.20972318031
$data[0,14] = 0.002925491708184224
$data[0,15] = 0.002925491708184224
$ws.Range("E9:T9").Value = $data

# Row 10: MuSCs -> ECs
$abcd = New-Object 'object[,]' 1,4
$abcd[0,0] = "MuSCs"
$abcd[0,1] = "Reln"
$abcd[0,2] = "Lrp8"
$abcd[0,3] = "ECs"
$ws.Range("A10:D10").Value = $abcd

$data = New-Object 'object[,]' 1,16
$data[0,0] = 3
$data[0,1] = 1
$data[0,2] = 4.62452
$data[0,3] = 13.87356
$data[0,4] = 0.9727340463997778
$data[0,5] = 0.9727340463997778
$data[0,6] = 3
$data[0,7] = 1
$data[0,8] = 0.4216986666666666
$data[0,9] = 1.265096
$data[0,10] = 0.2516921781598699
$data[0,11] = 0.2516921781598699
$data[0,12] = 1.950153917973333
$data[0,13] = 17.55138526176
$data[0,14] = 0.244829550908624
$data[0,15] = 0.244829550908624
$ws.Range("E10:T10").Value = $data

# Row 11: MuSCs -> FAPs
$abcd = New-Object 'object[,]' 1,4
$abcd[0,0] = "MuSCs"
$abcd[0,1] = "Reln"
$abcd[0,2] = "Lrp8"
$abcd[0,3] = "FAPs"
$ws.Range("A11:D11").Value = $abcd

$data = New-Object 'object[,]' 1,16
$data[0,0] = 3
$data[0,1] = 1
$data[0,2] = 4.62452
$data[0,3] = 13.87356
$data[0,4] = 0.9727340463997778
$data[0,5] = 0.9727340463997778
$data[0,6] = 3
$data[0,7] = 1
$data[0,8] = 0.437958
$data[0,9] = 1.313874
$data[0,10] = 0.2613966125002536
$data[0,11] = 0.2613966125002536
$data[0,12] = 2.02534553016
$data[0,13] = 18.22810977144
$data[0,14] = 0.2542693845925665
$data[0,15] = 0.2542693845925665
$ws.Range("E11:T11").Value = $data

# Row 12: MuSCs -> MuSCs
$abcd = New-Object 'object[,]' 1,4
$abcd[0,0] = "MuSCs"
$abcd[0,1] = "Reln"
$abcd[0,2] = "Lrp8"
$abcd[0,3] = "MuSCs"
$ws.Range("A12:D12").Value = $abcd

$data = New-Object 'object[,]' 1,16
$data[0,0] = 3
$data[0,1] = 1
$data[0,2] = 4.62452
$data[0,3] = 13.87356
$data[0,4] = 0.9727340463997778
$data[0,5] = 0.9727340463997778
$data[0,6] = 3
$data[0,7] = 1
$data[0,8] = 0.4328273333333333
$data[0,9] = 1.298482
$data[0,10] = 0.2583343579312433
$data[0,11] = 0.2583343579312433
$data[0,12] = 2.001618659546666
$data[0,13] = 18.01456793592
$data[0,14] = 0.2512906253145468
$data[0,15] = 0.2512906253145468
$ws.Range("E12:T12").Value = $data

# Row 13: MuSCs -> Resolving-Mac
$abcd = New-Object 'object[,]' 1,4
$abcd[0,0] = "MuSCs"
$abcd[0,1] = "Reln"
$abcd[0,2] = "Lrp8"
$abcd[0,3] = "Resolving-Mac"
$ws.Range("A13:D13").Value = $abcd

$data = New-Object 'object[,]' 1,16
$data[0,0] = 3
$data[0,1] = 1
$data[0,2] = 4.62452
$data[0,3] = 13.87356
$data[0,4] = 0.9727340463997778
$data[0,5] = 0.9727340463997778
$data[0,6] = 3
$data[0,7] = 1
$data[0,8] = 0.38297
$data[0,9] = 1.14891
$data[0,10] = 0.2285768514086331
$data[0,11] = 0.2285768514086331
$data[0,12] = 1.7710524244
$data[0,13] = 15.9394718196
$data[0,14] = 0.2223444855840404
$data[0,15] = 0.2223444855840404
$ws.Range("E13:T13").Value = $data

# Row 14: Resolving-Mac -> ECs
$abcd = New-Object 'object[,]' 1,4
$abcd[0,0] = "Resolving-Mac"
$abcd[0,1] = "Reln"
$abcd[0,2] = "Lrp8"
$abcd[0,3] = "ECs"
$ws.Range("A14:D14").Value = $abcd

$data = New-Object 'object[,]' 1,16
$data[0,0] = 1
$data[0,1] = 0.3333333333333333
$data[0,2] = 0.004225666666666667
$data[0,3] = 0.012677
$data[0,4] = 0.0008888381573446169
$data[0,5] = 0.000888838157344617
$data[0,6] = 3
$data[0,7] = 1
$data[0,8] = 0.4216986666666666
$data[0,9] = 1.265096
$data[0,10] = 0.2516921781598699
$data[0,11] = 0.2516921781598699
$data[0,12] = 0.001781957999111111
$data[0,13] = 0.016037621992
$data[0,14] = 0.0002237136118536718
$data[0,15] = 0.0002237136118536718
$ws.Range("E14:T14").Value = $data

# Row 15: Resolving-Mac -> FAPs
$abcd = New-Object 'object[,]' 1,4
$abcd[0,0] = "Resolving-Mac"
$abcd[0,1] = "Reln"
$abcd[0,2] = "Lrp8"
$abcd[0,3] = "FAPs"
$ws.Range("A15:D15").Value = $abcd

$data = New-Object 'object[,]' 1,16
$data[0,0] = 1
$data[0,1] = 0.3333333333333333
$data[0,2] = 0.004225666666666667
$data[0,3] = 0.012677
$data[0,4] = 0.0008888381573446169
$data[0,5] = 0.000888838157344617
$data[0,6] = 3
$data[0,7] = 1
$data[0,8] = 0.437958
$data[0,9] = 1.313874
$data[0,10] = 0.2613966125002536
$data[0,11] = 0.2613966125002536
$data[0,12] = 0.001850664522
$data[0,13] = 0.016655980698
$data[0,14] = 0.0002323392833908503
$data[0,15] = 0.0002323392833908503
$ws.Range("E15:T15").Value = $data

# Row 16: Resolving-Mac -> MuSCs
$abcd = New-Object 'object[,]' 1,4
$abcd[0,0] = "Resolving-Mac"
$abcd[0,1] = "Reln"
$abcd[0,2] = "Lrp8"
$abcd[0,3] = "MuSCs"
$ws.Range("A16:D16").Value = $abcd

$data = New-Object 'object[,]' 1,16
$data[0,0] = 1
$data[0,1] = 0.3333333333333333
$data[0,2] = 0.004225666666666667
$data[0,3] = 0.012677
$data[0,4] = 0.0008888381573446169
$data[0,5] = 0.000888838157344617
$data[0,6] = 3
$data[0,7] = 1
$data[0,8] = 0.4328273333333333
$data[0,9] = 1.298482
$data[0,10] = 0.2583343579312433
$data[0,11] = 0.2583343579312433
$data[0,12] = 0.001828984034888889
$data[0,13] = 0.016460856314
$data[0,14] = 0.000229617434682411
$data[0,15] = 0.0002296174346824111
$ws.Range("E16:T16").Value = $data

# Row 17: Resolving-Mac -> Resolving-Mac
$abcd = New-Object 'object[,]' 1,4
$abcd[0,0] = "Resolving-Mac"
$abcd[0,1] = "Reln"
$abcd[0,2] = "Lrp8"
$abcd[0,3] = "Resolving-Mac"
$ws.Range("A17:D17").Value = $abcd

$data = New-Object 'object[,]' 1,16
$data[0,0] = 1
$data[0,1] = 0.3333333333333333
$data[0,2] = 0.004225666666666667
$data[0,3] = 0.012677
$data[0,4] = 0.0008888381573446169
$data[0,5] = 0.000888838157344617
$data[0,6] = 3
$data[0,7] = 1
$data[0,8] = 0.38297
$data[0,9] = 1.14891
$data[0,10] = 0.2285768514086331
$data[0,11] = 0.2285768514086331
$data[0,12] = 0.001618303563333333
$data[0,13] = 0.01456473207
$data[0,14] = 0.0002031678274176837
$data[0,15] = 0.0002031678274176838
$ws.Range("E17:T17").Value = $data
